$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension-driving values: add columns P (16th) and Q (17th)
$ws.Cells.Item(1, 16).Value = 14   # P1
$ws.Cells.Item(1, 17).Value = 15   # Q1

# Copy header style (bold, bordered, centered) from O1 onto the new header cells
$ws.Cells.Item(1, 15).Copy() | Out-Null
$ws.Range($ws.Cells.Item(1, 16), $ws.Cells.Item(1, 17)).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 16).Value = 14   # restore value after paste-format
$ws.Cells.Item(1, 17).Value = 15

# Row 2
$ws.Cells.Item(2, 2).Value = 3.457262170408171
$ws.Cells.Item(2, 3).Value = 1.007702210158641
$ws.Cells.Item(2, 4).Value = 0.3492156199738474
$ws.Cells.Item(2, 5).Value = 1.343180230842464
$ws.Cells.Item(2, 6).Value = 6.463794665469322
$ws.Cells.Item(2, 7).Value = 0.0007918837209675011
$ws.Cells.Item(2, 8).Value = 0.00999529574372202
$ws.Cells.Item(2, 9).Value = 0.002670933310996659
$ws.Cells.Item(2, 16).Value = 0   # P2
$ws.Cells.Item(2, 17).Value = 0   # Q2

# Row 3
$ws.Cells.Item(3, 2).Value = 2.997853015548628
$ws.Cells.Item(3, 3).Value = 0.869853693225707
$ws.Cells.Item(3, 4).Value = 0.3064966316520383
$ws.Cells.Item(3, 5).Value = 1.155751310872887
$ws.Cells.Item(3, 6).Value = 5.654258087503649
$ws.Cells.Item(3, 7).Value = 0.0007995267261627653
$ws.Cells.Item(3, 8).Value = 0.00566132105441941
$ws.Cells.Item(3, 9).Value = 0.0006838291904367466
$ws.Cells.Item(3, 16).Value = 0   # P3
$ws.Cells.Item(3, 17).Value = 0   # Q3

# Row 4
$ws.Cells.Item(4, 2).Value = 2.717878933285135
$ws.Cells.Item(4, 3).Value = 0.7870470786523356
$ws.Cells.Item(4, 4).Value = 0.2802437474734631
$ws.Cells.Item(4, 5).Value = 1.042418103720038
$ws.Cells.Item(4, 6).Value = 5.159655295851422
$ws.Cells.Item(4, 7).Value = 0.00080432863316101
$ws.Cells.Item(4, 8).Value = 0.003562661497822361
$ws.Cells.Item(4, 9).Value = 0.0004342247631967666
$ws.Cells.Item(4, 16).Value = 0   # P4
$ws.Cells.Item(4, 17).Value = 0   # Q4

# Row 5
$ws.Cells.Item(5, 2).Value = 2.604150420224698
$ws.Cells.Item(5, 3).Value = 0.7552167163809997
$ws.Cells.Item(5, 4).Value = 0.2683125171909211
$ws.Cells.Item(5, 5).Value = 0.9965441341777392
$ws.Cells.Item(5, 6).Value = 4.944718871251297
$ws.Cells.Item(5, 7).Value = 0.0008063365970838446
$ws.Cells.Item(5, 8).Value = 0.00283109922356406
$ws.Cells.Item(5, 9).Value = 0.000622890166188661
$ws.Cells.Item(5, 16).Value = 0   # P5
$ws.Cells.Item(5, 17).Value = 0   # Q5

# Row 6
$ws.Cells.Item(6, 2).Value = 2.585207879043878
$ws.Cells.Item(6, 3).Value = 0.7518079961090791
$ws.Cells.Item(6, 4).Value = 0.2648516623235935
$ws.Cells.Item(6, 5).Value = 0.9888719871438241
$ws.Cells.Item(6, 6).Value = 4.892255946361189
$ws.Cells.Item(6, 7).Value = 0.0008066988388415279
$ws.Cells.Item(6, 8).Value = 0.002711539273524188
$ws.Cells.Item(6, 9).Value = 0.0007586302123687716
$ws.Cells.Item(6, 16).Value = 0   # P6
$ws.Cells.Item(6, 17).Value = 0   # Q6

# Row 7
$ws.Cells.Item(7, 2).Value = 2.716100113317282
$ws.Cells.Item(7, 3).Value = 0.7916813407236987
$ws.Cells.Item(7, 4).Value = 0.2760110113271139
$ws.Cells.Item(7, 5).Value = 1.041580167724106
$ws.Cells.Item(7, 6).Value = 5.110605349627804
$ws.Cells.Item(7, 7).Value = 0.0008044293528601984
$ws.Cells.Item(7, 8).Value = 0.00353482695453966
$ws.Cells.Item(7, 9).Value = 0.0006621355425160402
$ws.Cells.Item(7, 16).Value = 0   # P7
$ws.Cells.Item(7, 17).Value = 0   # Q7

# Row 8
$ws.Cells.Item(8, 2).Value = 3.297930760598206
$ws.Cells.Item(8, 3).Value = 0.9664754913100921
$ws.Cells.Item(8, 4).Value = 0.3289524093292329
$ws.Cells.Item(8, 5).Value = 1.277802295241244
$ws.Cells.Item(8, 6).Value = 6.121670233424283
$ws.Cells.Item(8, 7).Value = 0.0007945947571169861
$ws.Cells.Item(8, 8).Value = 0.00834023203055928
$ws.Cells.Item(8, 9).Value = 0.001998973943032389
$ws.Cells.Item(8, 16).Value = 0   # P8
$ws.Cells.Item(8, 17).Value = 0   # Q8

# Row 9
$ws.Cells.Item(9, 2).Value = 4.461276162754871
$ws.Cells.Item(9, 3).Value = 1.317628976875312
$ws.Cells.Item(9, 4).Value = 0.4409129175473083
$ws.Cells.Item(9, 5).Value = 1.758791490145597
$ws.Cells.Item(9, 6).Value = 8.221729476999883
$ws.Cells.Item(9, 7).Value = 0.0007759650060087116
$ws.Cells.Item(9, 8).Value = 0.02293966575955908
$ws.Cells.Item(9, 9).Value = 0.01318788530029646
$ws.Cells.Item(9, 16).Value = 0   # P9
$ws.Cells.Item(9, 17).Value = 0   # Q9

# Row 10
$ws.Cells.Item(10, 2).Value = 5.332556541617578
$ws.Cells.Item(10, 3).Value = 1.589096750276838
$ws.Cells.Item(10, 4).Value = 0.5003877811189454
$ws.Cells.Item(10, 5).Value = 2.021280822810638
$ws.Cells.Item(10, 6).Value = 9.52940342855976
$ws.Cells.Item(10, 7).Value = 0.0007632755018147897
$ws.Cells.Item(10, 8).Value = 0.03654587955837574
$ws.Cells.Item(10, 9).Value = 0.02786243332301552
$ws.Cells.Item(10, 16).Value = 0   # P10
$ws.Cells.Item(10, 17).Value = 0   # Q10

# Row 11
$ws.Cells.Item(11, 2).Value = 5.696287717332439
$ws.Cells.Item(11, 3).Value = 1.68929629314573
$ws.Cells.Item(11, 4).Value = 0.3475534396946358
$ws.Cells.Item(11, 5).Value = 1.319623735495057
$ws.Cells.Item(11, 6).Value = 7.984121445527933
$ws.Cells.Item(11, 7).Value = 0.0007625219121048961
$ws.Cells.Item(11, 8).Value = 0.05102710592101189
$ws.Cells.Item(11, 9).Value = 0.02986244529849014
$ws.Cells.Item(11, 16).Value = 0   # P11
$ws.Cells.Item(11, 17).Value = 0   # Q11

# Row 12
$ws.Cells.Item(12, 2).Value = 5.821453138132654
$ws.Cells.Item(12, 3).Value = 1.71081905899149
$ws.Cells.Item(12, 4).Value = 0.2378295168972926
$ws.Cells.Item(12, 5).Value = 0.8095450099342685
$ws.Cells.Item(12, 6).Value = 6.621540193137264
$ws.Cells.Item(12, 7).Value = 0.0007639042150792652
$ws.Cells.Item(12, 8).Value = 0.08537059510579326
$ws.Cells.Item(12, 9).Value = 0.02834500312373667
$ws.Cells.Item(12, 16).Value = 0   # P12
$ws.Cells.Item(12, 17).Value = 0   # Q12

# Row 13
$ws.Cells.Item(13, 2).Value = 5.769062860697261
$ws.Cells.Item(13, 3).Value = 1.681431827745541
$ws.Cells.Item(13, 4).Value = 0.1489218976295632
$ws.Cells.Item(13, 5).Value = 0.4192333535295205
$ws.Cells.Item(13, 6).Value = 5.257492853800528
$ws.Cells.Item(13, 7).Value = 0.0007670939994285801
$ws.Cells.Item(13, 8).Value = 0.1363447470369579
$ws.Cells.Item(13, 9).Value = 0.02441139653021018
$ws.Cells.Item(13, 16).Value = 0   # P13
$ws.Cells.Item(13, 17).Value = 0   # Q13

# Row 14
$ws.Cells.Item(14, 2).Value = 5.653258874752623
$ws.Cells.Item(14, 3).Value = 1.640501665964337
$ws.Cells.Item(14, 4).Value = 0.09914364752189897
$ws.Cells.Item(14, 5).Value = 0.2200208786317432
$ws.Cells.Item(14, 6).Value = 4.320288365004302
$ws.Cells.Item(14, 7).Value = 0.0007700778533965794
$ws.Cells.Item(14, 8).Value = 0.18241339872894
$ws.Cells.Item(14, 9).Value = 0.02096641771609153
$ws.Cells.Item(14, 16).Value = 0   # P14
$ws.Cells.Item(14, 17).Value = 0   # Q14

# Row 15
$ws.Cells.Item(15, 2).Value = 5.587027398541977
$ws.Cells.Item(15, 3).Value = 1.621696695691526
$ws.Cells.Item(15, 4).Value = 0.08790478637317278
$ws.Cells.Item(15, 5).Value = 0.180148105521539
$ws.Cells.Item(15, 6).Value = 4.065077586257132
$ws.Cells.Item(15, 7).Value = 0.0007712711082461959
$ws.Cells.Item(15, 8).Value = 0.1938160907104987
$ws.Cells.Item(15, 9).Value = 0.01972158660848322
$ws.Cells.Item(15, 16).Value = 0   # P15
$ws.Cells.Item(15, 17).Value = 0   # Q15

# Row 16
$ws.Cells.Item(16, 2).Value = 5.224085017559617
$ws.Cells.Item(16, 3).Value = 1.516492032032318
$ws.Cells.Item(16, 4).Value = 0.08526566423042681
$ws.Cells.Item(16, 5).Value = 0.1700411051112454
$ws.Cells.Item(16, 6).Value = 3.828971836814645
$ws.Cells.Item(16, 7).Value = 0.00077591099244641
$ws.Cells.Item(16, 8).Value = 0.1772398316277872
$ws.Cells.Item(16, 9).Value = 0.01491254195543057
$ws.Cells.Item(16, 16).Value = 0   # P16
$ws.Cells.Item(16, 17).Value = 0   # Q16

# Row 17
$ws.Cells.Item(17, 2).Value = 5.007352928785792
$ws.Cells.Item(17, 3).Value = 1.457340512143503
$ws.Cells.Item(17, 4).Value = 0.1063005318635462
$ws.Cells.Item(17, 5).Value = 0.2527426306022065
$ws.Cells.Item(17, 6).Value = 4.119620197415884
$ws.Cells.Item(17, 7).Value = 0.0007779609381707517
$ws.Cells.Item(17, 8).Value = 0.1380845818779335
$ws.Cells.Item(17, 9).Value = 0.01298819220147251
$ws.Cells.Item(17, 16).Value = 0   # P17
$ws.Cells.Item(17, 17).Value = 0   # Q17

# Row 18
$ws.Cells.Item(18, 2).Value = 4.891181072940014
$ws.Cells.Item(18, 3).Value = 1.426784906721878
$ws.Cells.Item(18, 4).Value = 0.1630917184647416
$ws.Cells.Item(18, 5).Value = 0.488672800946965
$ws.Cells.Item(18, 6).Value = 4.989782952390868
$ws.Cells.Item(18, 7).Value = 0.0007777175351106891
$ws.Cells.Item(18, 8).Value = 0.08624958060428156
$ws.Cells.Item(18, 9).Value = 0.0128610273449894
$ws.Cells.Item(18, 16).Value = 0   # P18
$ws.Cells.Item(18, 17).Value = 0   # Q18

# Row 19
$ws.Cells.Item(19, 2).Value = 4.865410637195282
$ws.Cells.Item(19, 3).Value = 1.433989826580387
$ws.Cells.Item(19, 4).Value = 0.2585924094822332
$ws.Cells.Item(19, 5).Value = 0.9316599478342198
$ws.Cells.Item(19, 6).Value = 6.311405617291229
$ws.Cells.Item(19, 7).Value = 0.0007753471943131185
$ws.Cells.Item(19, 8).Value = 0.04486607482094485
$ws.Cells.Item(19, 9).Value = 0.01491509726212481
$ws.Cells.Item(19, 16).Value = 0   # P19
$ws.Cells.Item(19, 17).Value = 0   # Q19

# Row 20
$ws.Cells.Item(20, 2).Value = 5.100185226391659
$ws.Cells.Item(20, 3).Value = 1.532472938516833
$ws.Cells.Item(20, 4).Value = 0.4702190910753643
$ws.Cells.Item(20, 5).Value = 1.945769640609583
$ws.Cells.Item(20, 6).Value = 9.0203272447026
$ws.Cells.Item(20, 7).Value = 0.0007668416119520997
$ws.Cells.Item(20, 8).Value = 0.03244554102579755
$ws.Cells.Item(20, 9).Value = 0.02379408360536761
$ws.Cells.Item(20, 16).Value = 0   # P20
$ws.Cells.Item(20, 17).Value = 0   # Q20

# Row 21
$ws.Cells.Item(21, 2).Value = 5.784892976581773
$ws.Cells.Item(21, 3).Value = 1.750257957264921
$ws.Cells.Item(21, 4).Value = 0.5534414014092022
$ws.Cells.Item(21, 5).Value = 2.318371252717682
$ws.Cells.Item(21, 6).Value = 10.45960704905036
$ws.Cells.Item(21, 7).Value = 0.000756201365715177
$ws.Cells.Item(21, 8).Value = 0.0464444912688986
$ws.Cells.Item(21, 9).Value = 0.03888522512563775
$ws.Cells.Item(21, 16).Value = 0   # P21
$ws.Cells.Item(21, 17).Value = 0   # Q21

# Row 22
$ws.Cells.Item(22, 2).Value = 6.238143054835405
$ws.Cells.Item(22, 3).Value = 1.888477988335353
$ws.Cells.Item(22, 4).Value = 0.6034973891750326
$ws.Cells.Item(22, 5).Value = 2.513470465521181
$ws.Cells.Item(22, 6).Value = 11.35561038052936
$ws.Cells.Item(22, 7).Value = 0.0007494717359220136
$ws.Cells.Item(22, 8).Value = 0.05622256578920748
$ws.Cells.Item(22, 9).Value = 0.05032041316143676
$ws.Cells.Item(22, 16).Value = 0   # P22
$ws.Cells.Item(22, 17).Value = 0   # Q22

# Row 23
$ws.Cells.Item(23, 2).Value = 5.996341369730999
$ws.Cells.Item(23, 3).Value = 1.808243448853432
$ws.Cells.Item(23, 4).Value = 0.5820632899846885
$ws.Cells.Item(23, 5).Value = 2.409407682906746
$ws.Cells.Item(23, 6).Value = 10.93547678253219
$ws.Cells.Item(23, 7).Value = 0.0007529736297253577
$ws.Cells.Item(23, 8).Value = 0.05097258498324031
$ws.Cells.Item(23, 9).Value = 0.04393539907997379
$ws.Cells.Item(23, 16).Value = 0   # P23
$ws.Cells.Item(23, 17).Value = 0   # Q23

# Row 24
$ws.Cells.Item(24, 2).Value = 5.093396454990398
$ws.Cells.Item(24, 3).Value = 1.522929633599688
$ws.Cells.Item(24, 4).Value = 0.494095960497134
$ws.Cells.Item(24, 5).Value = 2.024357664267967
$ws.Cells.Item(24, 6).Value = 9.282577626360364
$ws.Cells.Item(24, 7).Value = 0.0007664096614360714
$ws.Cells.Item(24, 8).Value = 0.03328234609949288
$ws.Cells.Item(24, 9).Value = 0.02388743652415215
$ws.Cells.Item(24, 16).Value = 0   # P24
$ws.Cells.Item(24, 17).Value = 0   # Q24

# Row 25
$ws.Cells.Item(25, 2).Value = 4.142656613778001
$ws.Cells.Item(25, 3).Value = 1.229542174013773
$ws.Cells.Item(25, 4).Value = 0.4026513398619613
$ws.Cells.Item(25, 5).Value = 1.625711697670212
$ws.Cells.Item(25, 6).Value = 7.560761790821061
$ws.Cells.Item(25, 7).Value = 0.0007810655889403182
$ws.Cells.Item(25, 8).Value = 0.01826779157226266
$ws.Cells.Item(25, 9).Value = 0.009293847328576632
$ws.Cells.Item(25, 16).Value = 0   # P25
$ws.Cells.Item(25, 17).Value = 0   # Q25
